# Atualize o programa de criação de planilha de acordo com o modelo
#
# Rebuilds the "Dashboard" sheet from the old compact layout (cols B-F)
# into the new, wider layout (cols A-H) with a "Valor da Carteira" summary
# block, per the target model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Start clean: drop the old merges (they overlap the new layout) and
# wipe existing cell content/formatting before laying out the new sheet.
# ---------------------------------------------------------------------
$ws.Cells.UnMerge()
$ws.Cells.Clear()

# ---------------------------------------------------------------------
# Column widths: A:H -> 27 characters.
# (ColumnWidth adds Excel's standard ~5-pixel/0.8333-char padding on
# top of what gets stored in the OOXML, so back that out here to land
# on a stored width of exactly 27.)
# ---------------------------------------------------------------------
$ws.Range("A1:H1").ColumnWidth = 26.166666666666668

# ---------------------------------------------------------------------
# Row 1-2: big title, "Resumo da Carteira", merged across A1:H2
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Resumo da Carteira"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 20
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108

$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(2).RowHeight = 27

$ws.Range("A1:H2").Merge()

# ---------------------------------------------------------------------
# Row 3-4: section headers "Ações" / "Moedas"
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Ações"
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Size = 16
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4108

$ws.Range("E3").Value = "Moedas"
$ws.Range("E3").Font.Bold = $true
$ws.Range("E3").Font.Size = 16
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").VerticalAlignment = -4108

$ws.Range("A3:D4").Merge()
$ws.Range("E3:H4").Merge()

# ---------------------------------------------------------------------
# Row 5: column headers for the two tables
# ---------------------------------------------------------------------
$headerCells = "A5","B5","C5","D5","E5","F5","G5","H5"
$headerTexts = "Nome","Quantidade","Valor da ação (R$)","Valor acumulado (R$)","Nome","Quantidade","Valor da ação (R$)","Valor acumulado (R$)"
for ($i = 0; $i -lt $headerCells.Length; $i++) {
  $cell = $ws.Range($headerCells[$i])
  $cell.Value = $headerTexts[$i]
  $cell.Font.Bold = $true
  $cell.Font.Size = 12
  $cell.HorizontalAlignment = -4108
  $cell.VerticalAlignment = -4108
  $cell.Merge()
}

# ---------------------------------------------------------------------
# Rows 6-8: "Ações" data (Nome, Quantidade, Valor da ação, Valor acumulado)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "VALE3"
$ws.Range("B6").Value = 1000
$ws.Range("C6").Value = 84.25
$ws.Range("C6").NumberFormat = "R$#,##0.00"
$ws.Range("D6").Formula = "=B6*C6"
$ws.Range("D6").NumberFormat = "R$#,##0.00"

$ws.Range("A7").Value = "MGLU3"
$ws.Range("B7").Value = 1000
$ws.Range("C7").Value = 4.07
$ws.Range("C7").NumberFormat = "R$#,##0.00"
$ws.Range("D7").Formula = "=B7*C7"
$ws.Range("D7").NumberFormat = "R$#,##0.00"

$ws.Range("A8").Value = "ITUB4"
$ws.Range("B8").Value = 375
$ws.Range("C8").Value = 26.02
$ws.Range("C8").NumberFormat = "R$#,##0.00"
$ws.Range("D8").Formula = "=B8*C8"
$ws.Range("D8").NumberFormat = "R$#,##0.00"

# ---------------------------------------------------------------------
# Rows 6-7: "Moedas" data (Nome, Quantidade, Valor da ação, Valor acumulado)
# ---------------------------------------------------------------------
$ws.Range("E6").Value = "CAD"
$ws.Range("F6").Value = 150
$ws.Range("G6").Value = 3.74
$ws.Range("G6").NumberFormat = "R#,##0.00"
$ws.Range("H6").Formula = "=F6*G6"
$ws.Range("H6").NumberFormat = "R#,##0.00"

$ws.Range("E7").Value = "CHF"
$ws.Range("F7").Value = 500
$ws.Range("G7").Value = 4.97
$ws.Range("G7").NumberFormat = "R#,##0.00"
$ws.Range("H7").Formula = "=F7*G7"
$ws.Range("H7").NumberFormat = "R#,##0.00"

# ---------------------------------------------------------------------
# Row 9: "Total Ações" / "Total Moedas" subtotal row
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Total Ações"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Size = 12
$ws.Range("A9").HorizontalAlignment = -4108
$ws.Range("A9").VerticalAlignment = -4108
$ws.Range("A9").Merge()

$ws.Range("B9").Formula = "=SUM(B6:B8)"

$ws.Range("C9").Formula = "=SUM(C6:C8)"
$ws.Range("C9").NumberFormat = "R#,##0.00"

$ws.Range("D9").Formula = "=SUM(D6:D8)"
$ws.Range("D9").NumberFormat = "R#,##0.00"

$ws.Range("E9").Value = "Total Moedas"
$ws.Range("E9").Font.Bold = $true
$ws.Range("E9").Font.Size = 12
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").Merge()

$ws.Range("F9").Formula = "=SUM(F6:F8)"

$ws.Range("G9").Formula = "=SUM(G6:G8)"
$ws.Range("G9").NumberFormat = "R#,##0.00"

$ws.Range("H9").Formula = "=SUM(H6:H8)"
$ws.Range("H9").NumberFormat = "R#,##0.00"

# ---------------------------------------------------------------------
# Row 12-13: "Valor da Carteira" sub-title, merged D12:E13
# ---------------------------------------------------------------------
$ws.Range("D12").Value = "Valor da Carteira"
$ws.Range("D12").Font.Bold = $true
$ws.Range("D12").Font.Size = 16
$ws.Range("D12").HorizontalAlignment = -4108
$ws.Range("D12").VerticalAlignment = -4108
$ws.Range("D12:E13").Merge()

# ---------------------------------------------------------------------
# Row 14: headers for the grand-total block
# ---------------------------------------------------------------------
$ws.Range("D14").Value = "Quantidade"
$ws.Range("D14").Font.Bold = $true
$ws.Range("D14").Font.Size = 12
$ws.Range("D14").HorizontalAlignment = -4108
$ws.Range("D14").VerticalAlignment = -4108
$ws.Range("D14").Merge()

$ws.Range("E14").Value = "Valor acumulado total (R$)"
$ws.Range("E14").Font.Bold = $true
$ws.Range("E14").Font.Size = 12
$ws.Range("E14").HorizontalAlignment = -4108
$ws.Range("E14").VerticalAlignment = -4108
$ws.Range("E14").Merge()

# ---------------------------------------------------------------------
# Row 15: grand totals
# ---------------------------------------------------------------------
$ws.Range("D15").Formula = "=B9+F9"
$ws.Range("D15").NumberFormat = "R#,##0.00"

$ws.Range("E15").Formula = "=D9+H9"
$ws.Range("E15").NumberFormat = "R#,##0.00"
